# "day 6 apache poi" - update FirstName values for rows 2 and 3 of the
# Employees sheet (Steven -> Tom, MadamM -> Adam), and move the active
# selection to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Tom"
$ws.Range("A3").Value = "Adam"
$ws.Range("A3").Select()
